$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (AD1, AE1, AF1), copying the existing header style
# from A1 so the new cells match the other header cells (bold, centered,
# bordered).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2-43): every row gets the
# same team record (86 wins, 76 losses, 0 ties).
$ws.Range("AD2:AD43").Value = 86
$ws.Range("AE2:AE43").Value = 76
$ws.Range("AF2:AF43").Value = 0

Write-Output "done"
